$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, date range) ---
# "Volume 32   Number  24" -> "...25"  (chars 21-22 are "24")
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "25"

# "Report Covering the Week  6/9/2025  Through  6/15/2025"
#  -> "...6/16/2025  Through  6/22/2025"
$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "6/16/2025"
$c9.Characters(47, 9).Text = "6/22/2025"

# --- Cells changing type/style (copy format+value from a pristine same-style source cell) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("F14").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 2
$ws.Range("F14").Copy($ws.Range("C25"))
$ws.Range("C25").Value = 2
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("F14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))

# --- Simple numeric value updates ---
$ws.Range("G15").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -53.846153846153
$ws.Range("I16").Value = 46
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = 9.523809523809
$ws.Range("L16").Value = -30.303030303030
$ws.Range("M16").Value = -44.578313253012
$ws.Range("N16").Value = -84.246575342465
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 23.529411764705
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = -0.961538461538
$ws.Range("L17").Value = 1.980198019801
$ws.Range("M17").Value = 49.275362318840
$ws.Range("N17").Value = -8.035714285714
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = 8.333333333333
$ws.Range("L18").Value = -7.142857142857
$ws.Range("M18").Value = -78.151260504201
$ws.Range("N18").Value = -91.333333333333
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("I19").Value = 120
$ws.Range("J19").Value = 153
$ws.Range("K19").Value = -21.568627450980
$ws.Range("L19").Value = 10.091743119266
$ws.Range("M19").Value = 16.504854368932
$ws.Range("N19").Value = -33.701657458563
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 22.222222222222
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 66
$ws.Range("K20").Value = -28.787878787878
$ws.Range("L20").Value = -38.961038961039
$ws.Range("M20").Value = -7.843137254901
$ws.Range("N20").Value = -95.378564405113
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -36.842105263157
$ws.Range("F21").Value = 58
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -13.432835820895
$ws.Range("I21").Value = 350
$ws.Range("J21").Value = 397
$ws.Range("K21").Value = -11.838790931989
$ws.Range("L21").Value = -10.256410256410
$ws.Range("M21").Value = -19.540229885057
$ws.Range("N21").Value = -81.742305685967
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -44.444444444444
$ws.Range("I23").Value = 51
$ws.Range("J23").Value = 58
$ws.Range("K23").Value = -12.068965517241
$ws.Range("L23").Value = -13.559322033898
$ws.Range("M23").Value = 142.857142857143
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 8
$ws.Range("F24").Value = 28
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = -41.666666666666
$ws.Range("I24").Value = 237
$ws.Range("J24").Value = 258
$ws.Range("K24").Value = -8.139534883720
$ws.Range("L24").Value = -17.132867132867
$ws.Range("M24").Value = -1.25
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = -80
$ws.Range("I25").Value = 54
$ws.Range("J25").Value = 52
$ws.Range("K25").Value = 3.846153846153
$ws.Range("L25").Value = -29.870129870129
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -11.538461538461
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 153
$ws.Range("K26").Value = -3.921568627450
$ws.Range("L26").Value = -15.028901734104
$ws.Range("M26").Value = -34.666666666666
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 16
$ws.Range("K28").Value = 6.666666666666
$ws.Range("L28").Value = 45.454545454545
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
